# fix error data of grouptc-cuckoo
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated GroupTC-HS timings (column E) and recalculated speedups (column I)
$ws.Range("E2").Value = 4.262
$ws.Range("I2").Value = 2.70178320037541

$ws.Range("E3").Value = 12.015
$ws.Range("I3").Value = 2.975280898876405

$ws.Range("E4").Value = 34.637
$ws.Range("I4").Value = 3.062620896728931

$ws.Range("E5").Value = 94.834
$ws.Range("I5").Value = 3.17860682877449

$ws.Range("E6").Value = 259.004
$ws.Range("I6").Value = 3.180997204676376

$ws.Range("E7").Value = 678.495
$ws.Range("I7").Value = 3.244557439627411

$ws.Range("E8").Value = 1761.689
$ws.Range("I8").Value = 3.278232423543542

$ws.Range("E9").Value = 4579.446
$ws.Range("I9").Value = 3.276141699236109
